$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values that look numeric (single decimal point) need to be
# forced to text so Excel does not auto-convert them to floating point numbers,
# matching the source data which stores them as literal text strings.
function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

$ws.Range("D2").Value = "28.005.65"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.882.77"
$ws.Range("E3").Value = "  +0.83%  "
Set-TextValue $ws.Range("D4") "1.010"
$ws.Range("E4").Value = "  +0.62%  "
Set-TextValue $ws.Range("D5") "335.93"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +1.45%  "
Set-TextValue $ws.Range("D8") "0.3953"
$ws.Range("E8").Value = "  +0.60%  "
Set-TextValue $ws.Range("D9") "47.32"
$ws.Range("E9").Value = "  -1.21%  "
Set-TextValue $ws.Range("D10") "0.08037"
$ws.Range("E10").Value = "  +0.00%  "
Set-TextValue $ws.Range("D11") "1.021"
$ws.Range("E11").Value = "  -0.19%  "
Set-TextValue $ws.Range("D12") "21.98"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "1.885.75"
$ws.Range("E13").Value = "  +3.12%  "
Set-TextValue $ws.Range("D14") "6.063"
$ws.Range("E14").Value = "  +1.98%  "
Set-TextValue $ws.Range("D15") "7.223"
$ws.Range("E15").Value = "  +1.26%  "
Set-TextValue $ws.Range("D16") "1.013"
$ws.Range("E16").Value = "  +0.96%  "
Set-TextValue $ws.Range("D17") "88.93"
$ws.Range("E17").Value = "  +2.59%  "
Set-TextValue $ws.Range("D18") "0.06751"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "27.990.70"
$ws.Range("E22").Value = "  +1.20%  "
Set-TextValue $ws.Range("D23") "5.523"
$ws.Range("E23").Value = "  +0.53%  "
Set-TextValue $ws.Range("D24") "11.03"
$ws.Range("E24").Value = "  +0.70%  "
Set-TextValue $ws.Range("D25") "2.344"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").Value = "2.107.32"
$ws.Range("E26").Value = "  +2.62%  "
Set-TextValue $ws.Range("D27") "158.90"
$ws.Range("E27").Value = "  +0.06%  "
Set-TextValue $ws.Range("D28") "19.93"
$ws.Range("E28").Value = "  -1.11%  "
Set-TextValue $ws.Range("D29") "2.113"
$ws.Range("E29").Value = "  +1.19%  "
Set-TextValue $ws.Range("D30") "5.513"
$ws.Range("E30").Value = "  -0.66%  "
Set-TextValue $ws.Range("D31") "121.73"
$ws.Range("E31").Value = "  -0.38%  "
Set-TextValue $ws.Range("D32") "0.9828"
$ws.Range("E32").Value = "  +1.69%  "
Set-TextValue $ws.Range("D33") "0.09582"
Set-TextValue $ws.Range("D34") "3.631"
$ws.Range("E34").Value = "  +1.05%  "
Set-TextValue $ws.Range("D35") "5.350"
$ws.Range("E35").Value = "  +0.54%  "
Set-TextValue $ws.Range("D36") "1.362"
$ws.Range("E36").Value = "  -5.84%  "
Set-TextValue $ws.Range("D37") "0.06088"
$ws.Range("E37").Value = "  -0.18%  "
Set-TextValue $ws.Range("D38") "0.02251"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  -1.46%  "
Set-TextValue $ws.Range("D40") "8.214"
$ws.Range("E40").Value = "  +1.37%  "
Set-TextValue $ws.Range("D41") "1.010"
$ws.Range("E41").Value = "  +0.74%  "
Set-TextValue $ws.Range("D42") "0.6009"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +0.40%  "
Set-TextValue $ws.Range("D44") "10.39"
$ws.Range("E44").Value = "  +1.44%  "
Set-TextValue $ws.Range("D45") "1.262"
$ws.Range("E45").Value = "  +0.28%  "
Set-TextValue $ws.Range("D46") "0.5694"
$ws.Range("E46").Value = "  -0.20%  "
Set-TextValue $ws.Range("D47") "12.21"
$ws.Range("E47").Value = "  -0.16%  "
Set-TextValue $ws.Range("D48") "1.939"
$ws.Range("E48").Value = "  +0.23%  "
Set-TextValue $ws.Range("D49") "3.347"
$ws.Range("E49").Value = "  -1.10%  "
Set-TextValue $ws.Range("D51") "112.93"
$ws.Range("E51").Value = "  -1.22%  "
